$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2,3,4,5,6,9,10,11,14,15)  # B,C,D,E,F,I,J,K,N,O

$data = @(
  @(0.3831433893469125, 0.05921439310179721, 0.2021402031676445, 0.1699234821612237, 1.364870054439962, 0.6573968767895622, 0.1860676919038298, 0.415350064829056, 1.534556796755091, 3.250231393374804),
  @(0.3465560993875272, 0.05264474798235597, 0.1959791953402572, 0.1656344628849737, 1.366129896015003, 0.662719455878154, 0.1820667567138514, 0.3741560588012192, 1.550230229227868, 3.265215660861841),
  @(0.3241472087435397, 0.04861178352794582, 0.192281238154635, 0.1630856878434663, 1.367601088904379, 0.6663089519006178, 0.1797137774834923, 0.3489117275351248, 1.560356193174051, 3.276085227157054),
  @(0.3150300246387019, 0.04696859921637042, 0.1907957639144371, 0.1620683984033526, 1.368376165317528, 0.6678525299136631, 0.1787810027246977, 0.3386373740623014, 1.564608957767437, 3.280934500691316),
  @(0.3135170231807365, 0.04669576890003668, 0.1905504024505404, 0.1619007699043777, 1.368515472839874, 0.6681137219910909, 0.1786276924287975, 0.3369321234719109, 1.565322757961491, 3.281765079315562),
  @(0.3240241912481281, 0.04858962173145187, 0.1922611174283162, 0.1630718817703922, 1.367610830882128, 0.6663294418489798, 0.1797010921068676, 0.3487731109090646, 1.560413035827869, 3.27614892610562),
  @(0.3705168334956852, 0.05694905091368696, 0.1999983296303043, 0.1684270822682592, 1.3651597232532, 0.659165397797068, 0.1846666810162034, 0.4011365619078617, 1.53985673146763, 3.255051596255242),
  @(0.4621111582371498, 0.07334594469617173, 0.2158405835517527, 0.1795988105055173, 1.365884641427115, 0.6476667727050653, 0.1952257254456242, 0.5041879994405747, 1.503531009705364, 3.226921202591512),
  @(0.5296417641425535, 0.0853930995506289, 0.2278837548161476, 0.1882137977710556, 1.369786322194997, 0.6407732323686801, 0.2034846020251138, 0.5801020371216907, 1.479269346664783, 3.214325847329945),
  @(0.5604103722753564, 0.09087335768853677, 0.2334493591923348, 0.1922211330145061, 1.372292362265995, 0.6379747378997678, 0.2073507274518249, 0.6146770701681419, 1.468758513699443, 3.210348869711112),
  @(0.5720681318272511, 0.09294852142048171, 0.2355693324520587, 0.1937512639850851, 1.373346389085171, 0.6369635425223699, 0.208830405964477, 0.6277751624272128, 1.464853953951687, 3.209094898844455),
  @(0.5695571521413854, 0.09250160316148026, 0.2351122088107473, 0.1934211617279828, 1.373114715202135, 0.6371791627589971, 0.2085110347021129, 0.6249540298429963, 1.465691504918935, 3.209353755284241),
  @(0.5613693394224697, 0.09104408496011729, 0.2336235228832351, 0.1923467648424051, 1.372376972992726, 0.6378905734765432, 0.207472147807195, 0.6157545561658253, 1.468435766561035, 3.21024065391228),
  @(0.556354876594213, 0.09015129863738025, 0.2327132708263804, 0.1916903099842813, 1.371938761593668, 0.6383326540356435, 0.2068378381537741, 0.6101202803372132, 1.470126560991538, 3.210816725511506),
  @(0.5276318870482157, 0.08503494359362662, 0.2275217715032625, 0.1879536811858387, 1.369637251080519, 0.6409629113745332, 0.2032341348127034, 0.5778432523850938, 1.479966840348682, 3.214621020069899),
  @(0.5100232599679657, 0.0818961525515931, 0.2243591750663825, 0.185683960796851, 1.368412546196552, 0.6426629184692771, 0.2010513063256667, 0.5580524574237131, 1.486138259767502, 3.217403733728759),
  @(0.4998998553976719, 0.08009080296342574, 0.2225483395651935, 0.1843867958557368, 1.367776952909892, 0.6436724754995673, 0.1998060739735337, 0.5466732483195358, 1.48973739654974, 3.219169251472664),
  @(0.4964730559453585, 0.07947954515536537, 0.221936635055684, 0.1839490286439016, 1.36757357548899, 0.6440197481799288, 0.1993862249973688, 0.5428211394769562, 1.490964502906575, 3.219795360381511),
  @(0.5118972556392407, 0.08223028282016287, 0.2246949904443341, 0.1859247159554513, 1.368535795716113, 0.6424786628826098, 0.2012826091638971, 0.560158819159625, 1.485476178347271, 3.217090434967503),
  @(0.5637741310534921, 0.09147219652447802, 0.2340604505088208, 0.1926619987016736, 1.37259081565594, 0.6376802975974343, 0.2077768692911377, 0.6184565266870266, 1.467627656563643, 3.209973310532632),
  @(0.5977154184745643, 0.09751175168457848, 0.240253536143598, 0.1971388493435384, 1.375853237375324, 0.6348271869941762, 0.2121125117478755, 0.656587899714367, 1.45640366638146, 3.206790841747477),
  @(0.5795971610612867, 0.09428840715474962, 0.2369416037177388, 0.1947427518873752, 1.374056033488856, 0.6363240558295686, 0.2097901570262621, 0.6362339033012745, 1.462353751773636, 3.208354977809051),
  @(0.5110500220557697, 0.08207922506305465, 0.2245431453047928, 0.1858158464737869, 1.368479861209806, 0.6425618645320554, 0.2011780069197187, 0.5592065367606551, 1.485775346167769, 3.217231561193273),
  @(0.4372894424898277, 0.06890993377530208, 0.2114835960723127, 0.1765049615753043, 1.365096661230716, 0.6505044641944728, 0.1922812386366246, 0.4762728017315396, 1.512931723578497, 3.233113664549563)
)

for ($r = 0; $r -lt $data.Count; $r++) {
  $rowNum = $r + 2
  $rowVals = $data[$r]
  for ($c = 0; $c -lt $cols.Count; $c++) {
    $ws.Cells.Item($rowNum, $cols[$c]).Value = $rowVals[$c]
  }
}